$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text/string updates (safe from numeric auto-conversion) ---
$ws.Range("D2").Value = "29.451.56"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.876.90"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.863.05"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "29.467.79"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E22").Value = "  +7.16%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E24").Value = "  +11.17%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.241.54"
$ws.Range("E38").Value = "  +8.35%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E41").Value = "  +5.47%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("E47").Value = "  +12.43%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E51").Value = "  +1.84%  "

# --- Numeric-looking text values that must be forced to remain text ---
# (Excel would otherwise auto-convert these to actual numbers)
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "0.7114"
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "241.96"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.07826"
$r.Style = "Normal"
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "25.14"
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.7271"
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "5.267"
$r.Style = "Normal"
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "90.96"
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "5.918"
$r.Style = "Normal"
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "247.48"
$r.Style = "Normal"
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "0.000007868"
$r.Style = "Normal"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "13.27"
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "7.975"
$r.Style = "Normal"
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "1.001"
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "0.1572"
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "163.54"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "9.014"
$r.Style = "Normal"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "18.29"
$r.Style = "Normal"
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "1.365"
$r.Style = "Normal"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.496"
$r.Style = "Normal"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "4.379"
$r.Style = "Normal"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "4.127"
$r.Style = "Normal"
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.05308"
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.931"
$r.Style = "Normal"
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "1.200"
$r.Style = "Normal"
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.7233"
$r.Style = "Normal"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "2.681"
$r.Style = "Normal"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.01863"
$r.Style = "Normal"
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "2.725"
$r.Style = "Normal"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.9078"
$r.Style = "Normal"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "74.05"
$r.Style = "Normal"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "6.149"
$r.Style = "Normal"
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "1.000"
$r.Style = "Normal"
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "103.31"
$r.Style = "Normal"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "0.5322"
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "1.773"
$r.Style = "Normal"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "2.915"
$r.Style = "Normal"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.00000000120"
$r.Style = "Normal"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.4319"
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "9.254"
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "7.070"
$r.Style = "Normal"
